$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I2 used to hold the text "20,50" (a shared string, since it used a comma as
# decimal separator). Replace it with the real numeric value 20.5 so it is
# stored as a number instead of text.
$ws.Cells.Item(2, 9).Value = 20.5

# Add a second line item ("Waffles") to the same invoice as row 4
# (2021-000003 / Jorge Cardano), showing several products can belong to one
# invoice.
$ws.Cells.Item(5, 1).Value = "2021-000003"
$ws.Cells.Item(5, 2).Style = $ws.Cells.Item(3, 2).Style
$ws.Cells.Item(5, 7).Value = "Waffles"
$ws.Cells.Item(5, 8).Value = 250
$ws.Cells.Item(5, 9).Value = 2.4
